# Increase font sizes throughout the resume, per the commit:
#   16pt -> 18pt  (name header)
#    9pt -> 10pt  (contact info / job dates / body text / bullets)
#   12pt -> 13pt  (section headers)
#   10pt -> 11pt  (overview / section intro paragraphs)
#   11pt -> 12pt  (job titles)
#
# Every run within a given paragraph shares the same size in this
# document, so walking paragraph-by-paragraph and rewriting
# Range.Font.Size is sufficient and keeps bold/italic/color intact.

$d = $word.ActiveDocument

$sizeMap = @{
    16 = 18
    9  = 10
    12 = 13
    10 = 11
    11 = 12
}

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $pr = $p.Range
    # Exclude the trailing paragraph-mark character so Word doesn't
    # stamp an explicit <w:rPr> onto <w:pPr> (which the source document
    # never had).
    $r = $d.Range($pr.Start, $pr.End - 1)
    $cur = $r.Font.Size
    if ($sizeMap.ContainsKey($cur)) {
        $r.Font.Size = $sizeMap[$cur]
    }
}
